$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from existing header cell (E1) to new headers so they match formatting
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-set header values after paste (paste special with formats only shouldn't touch values, but ensure correctness)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean outlier flag data for rows 2-17 across columns F (KNN), G (SVM), H (RF)
$values = @(
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $true,  $true),
    @($true,  $true,  $true),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $true),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
